# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker/period rows (B16:G47) are re-sorted by period and the new
# worker "PAOLA TATIANA OLIVER PEÑARANDA" (52501544) gains four extra,
# older periods (1808, 1810, 1811, 1812) ahead of the shared 1902-1908
# block that is now common to all four workers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=16; B='CC'; C='52501544'; D='PAOLA TATIANA OLIVER PEÑARANDA'; E='1808'; F=94889; G=5494911 },
    @{ Row=17; B='CC'; C='52501544'; D='PAOLA TATIANA OLIVER PEÑARANDA'; E='1810'; F=219797; G=5494911 },
    @{ Row=18; B='CC'; C='52501544'; D='PAOLA TATIANA OLIVER PEÑARANDA'; E='1811'; F=219797; G=5494911 },
    @{ Row=19; B='CC'; C='52501544'; D='PAOLA TATIANA OLIVER PEÑARANDA'; E='1812'; F=219797; G=5494911 },
    @{ Row=20; B='CC'; C='1033375418'; D='GEIDER MANUEL PEREZ OYOLA'; E='1902'; F=33125; G=828116 },
    @{ Row=21; B='CC'; C='1052735732'; D='GARIBALDIS GUERRERO FLOREZ'; E='1902'; F=33125; G=828116 },
    @{ Row=22; B='CC'; C='1052731447'; D='JOSE MANUEL ALVAREZ DIAZ'; E='1902'; F=33125; G=828116 },
    @{ Row=23; B='CC'; C='52501544'; D='PAOLA TATIANA OLIVER PEÑARANDA'; E='1902'; F=219797; G=5494911 },
    @{ Row=24; B='CC'; C='1033375418'; D='GEIDER MANUEL PEREZ OYOLA'; E='1903'; F=33125; G=828116 },
    @{ Row=25; B='CC'; C='1052735732'; D='GARIBALDIS GUERRERO FLOREZ'; E='1903'; F=33125; G=828116 },
    @{ Row=26; B='CC'; C='1052731447'; D='JOSE MANUEL ALVAREZ DIAZ'; E='1903'; F=33125; G=828116 },
    @{ Row=27; B='CC'; C='52501544'; D='PAOLA TATIANA OLIVER PEÑARANDA'; E='1903'; F=219797; G=5494911 },
    @{ Row=28; B='CC'; C='1033375418'; D='GEIDER MANUEL PEREZ OYOLA'; E='1904'; F=33125; G=828116 },
    @{ Row=29; B='CC'; C='1052735732'; D='GARIBALDIS GUERRERO FLOREZ'; E='1904'; F=33125; G=828116 },
    @{ Row=30; B='CC'; C='1052731447'; D='JOSE MANUEL ALVAREZ DIAZ'; E='1904'; F=33125; G=828116 },
    @{ Row=31; B='CC'; C='52501544'; D='PAOLA TATIANA OLIVER PEÑARANDA'; E='1904'; F=219797; G=5494911 },
    @{ Row=32; B='CC'; C='1033375418'; D='GEIDER MANUEL PEREZ OYOLA'; E='1905'; F=33125; G=828116 },
    @{ Row=33; B='CC'; C='1052735732'; D='GARIBALDIS GUERRERO FLOREZ'; E='1905'; F=33125; G=828116 },
    @{ Row=34; B='CC'; C='1052731447'; D='JOSE MANUEL ALVAREZ DIAZ'; E='1905'; F=33125; G=828116 },
    @{ Row=35; B='CC'; C='52501544'; D='PAOLA TATIANA OLIVER PEÑARANDA'; E='1905'; F=219797; G=5494911 },
    @{ Row=36; B='CC'; C='1033375418'; D='GEIDER MANUEL PEREZ OYOLA'; E='1906'; F=33125; G=828116 },
    @{ Row=37; B='CC'; C='1052735732'; D='GARIBALDIS GUERRERO FLOREZ'; E='1906'; F=33125; G=828116 },
    @{ Row=38; B='CC'; C='1052731447'; D='JOSE MANUEL ALVAREZ DIAZ'; E='1906'; F=33125; G=828116 },
    @{ Row=39; B='CC'; C='52501544'; D='PAOLA TATIANA OLIVER PEÑARANDA'; E='1906'; F=219797; G=5494911 },
    @{ Row=40; B='CC'; C='1033375418'; D='GEIDER MANUEL PEREZ OYOLA'; E='1907'; F=33125; G=828116 },
    @{ Row=41; B='CC'; C='1052735732'; D='GARIBALDIS GUERRERO FLOREZ'; E='1907'; F=33125; G=828116 },
    @{ Row=42; B='CC'; C='1052731447'; D='JOSE MANUEL ALVAREZ DIAZ'; E='1907'; F=33125; G=828116 },
    @{ Row=43; B='CC'; C='52501544'; D='PAOLA TATIANA OLIVER PEÑARANDA'; E='1907'; F=219797; G=5494911 },
    @{ Row=44; B='CC'; C='1033375418'; D='GEIDER MANUEL PEREZ OYOLA'; E='1908'; F=20979; G=828116 },
    @{ Row=45; B='CC'; C='1052735732'; D='GARIBALDIS GUERRERO FLOREZ'; E='1908'; F=20979; G=828116 },
    @{ Row=46; B='CC'; C='1052731447'; D='JOSE MANUEL ALVAREZ DIAZ'; E='1908'; F=20979; G=828116 },
    @{ Row=47; B='CC'; C='52501544'; D='PAOLA TATIANA OLIVER PEÑARANDA'; E='1908'; F=139205; G=5494911 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
}
